# Daily attendance processing - 2026-01-25 05:44:36
# Swap the order of "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# in column G ("Recorded By") wherever it occurs on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldVal = "dnasr281@gmail.com, System"
$newVal = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldVal) {
        $cell.Value = $newVal
    }
}
